# Weekly update: a new "Apio" price record (week of 2022-06-02) was added
# to the Vega Central Mapocho de Santiago sheet. Insert two new data rows
# right above the current row 270 (pushing the existing rows 270-287 down
# to 272-289) and populate them with the new weekly observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 270:271, shifting everything below down by two.
$ws.Range("A270:A271").EntireRow.Insert()

# New row 270: "Primera" quality record for 2022-06-02 (serial 44714)
$ws.Range("A270").Value = 9
$ws.Range("B270").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44714
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 100112017
$ws.Range("G270").Value = "Apio"
$ws.Range("H270").Value = "Americana (o)"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 79
$ws.Range("K270").Value = 7000
$ws.Range("L270").Value = 8000
$ws.Range("M270").Value = 7494
$ws.Range("N270").Value = "`$/docena de matas"
$ws.Range("O270").Value = "Región de Coquimbo"
$ws.Range("P270").Value = 1249
$ws.Range("Q270").Value = 6
$ws.Range("R270").Value = "Hortaliza"

# New row 271: "Segunda" quality record for 2022-06-02 (serial 44714)
$ws.Range("A271").Value = 9
$ws.Range("B271").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C271").Value = "Metropolitana"
$ws.Range("D271").Value = 44714
$ws.Range("E271").Value = 13
$ws.Range("F271").Value = 100112017
$ws.Range("G271").Value = "Apio"
$ws.Range("H271").Value = "Americana (o)"
$ws.Range("I271").Value = "Segunda"
$ws.Range("J271").Value = 34
$ws.Range("K271").Value = 6000
$ws.Range("L271").Value = 6000
$ws.Range("M271").Value = 6000
$ws.Range("N271").Value = "`$/docena de matas"
$ws.Range("O271").Value = "Región de Coquimbo"
$ws.Range("P271").Value = 1000
$ws.Range("Q271").Value = 6
$ws.Range("R271").Value = "Hortaliza"
